$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")  # row 33
$ws.Range("H33").Value = 668.3077
$ws.Range("I33").Value = 385.375
$ws.Range("K33").Value = 385.375
$ws.Range("M33").Value = -156.375

$ws = $wb.Worksheets.Item("ALC")  # row 86
$ws.Range("H86").Value = 333448060
$ws.Range("I86").Value = 500002500
$ws.Range("K86").Value = 500002500
$ws.Range("M86").Value = -500001377

$ws = $wb.Worksheets.Item("ALC")  # row 89
$ws.Range("H89").Value = 333448060
$ws.Range("I89").Value = 500002500
$ws.Range("K89").Value = 2500012500
$ws.Range("M89").Value = -2500006884

$ws = $wb.Worksheets.Item("ALC")  # row 132
$ws.Range("H132").Value = 2244.8845
$ws.Range("I132").Value = 1714.174
$ws.Range("K132").Value = 5142.522
$ws.Range("M132").Value = -2612.522

$ws = $wb.Worksheets.Item("ARM")  # row 2
$ws.Range("H2").Value = 1892.1538
$ws.Range("J2").Value = 2936.1428
$ws.Range("L2").Value = 2936.1428
$ws.Range("N2").Value = -3162.1428

$ws = $wb.Worksheets.Item("ARM")  # row 45
$ws.Range("H45").Value = 4561.0625
$ws.Range("J45").Value = 4192.8184
$ws.Range("L45").Value = 4192.8184
$ws.Range("N45").Value = -4946.8184

$ws = $wb.Worksheets.Item("ARM")  # row 63
$ws.Range("H63").Value = 3966.125
$ws.Range("I63").Value = 2432.25
$ws.Range("J63").Value = 5500
$ws.Range("K63").Value = 2432.25
$ws.Range("L63").Value = 5500
$ws.Range("M63").Value = -1746.25
$ws.Range("N63").Value = -6872

$ws = $wb.Worksheets.Item("ARM")  # row 66
$ws.Range("H66").Value = 3966.125
$ws.Range("I66").Value = 2432.25
$ws.Range("J66").Value = 5500
$ws.Range("K66").Value = 12161.25
$ws.Range("L66").Value = 27500
$ws.Range("M66").Value = -8729.25
$ws.Range("N66").Value = -34364

$ws = $wb.Worksheets.Item("ARM")  # row 74
$ws.Range("H74").Value = 2144.3928
$ws.Range("I74").Value = 2264.5264
$ws.Range("K74").Value = 2264.5264
$ws.Range("M74").Value = -1390.5264

$ws = $wb.Worksheets.Item("ARM")  # row 77
$ws.Range("H77").Value = 2144.3928
$ws.Range("I77").Value = 2264.5264
$ws.Range("K77").Value = 11322.632
$ws.Range("M77").Value = -6954.632000000001

$ws = $wb.Worksheets.Item("ARM")  # row 97
$ws.Range("H97").Value = 1588.64
$ws.Range("I97").Value = 1135.8
$ws.Range("J97").Value = 3400
$ws.Range("K97").Value = 1135.8
$ws.Range("L97").Value = 3400
$ws.Range("M97").Value = -639.8
$ws.Range("N97").Value = -4392

$ws = $wb.Worksheets.Item("ARM")  # row 102
$ws.Range("H102").Value = 2112.125
$ws.Range("J102").Value = 2345
$ws.Range("L102").Value = 2345
$ws.Range("N102").Value = -5589

$ws = $wb.Worksheets.Item("ARM")  # row 116
$ws.Range("H116").Value = 1892.1538
$ws.Range("J116").Value = 2936.1428
$ws.Range("L116").Value = 2936.1428
$ws.Range("N116").Value = -7524.1428

$ws = $wb.Worksheets.Item("ARM")  # row 122
$ws.Range("H122").Value = 5258.0557
$ws.Range("I122").Value = 4719.619
$ws.Range("K122").Value = 14158.857
$ws.Range("M122").Value = -11708.857

$ws = $wb.Worksheets.Item("ARM")  # row 132
$ws.Range("H132").Value = 224739.86
$ws.Range("I132").Value = 229801.81
$ws.Range("K132").Value = 689405.4299999999
$ws.Range("M132").Value = -686875.4299999999

$ws = $wb.Worksheets.Item("BSM")  # row 3
$ws.Range("H3").Value = 1892.1538
$ws.Range("J3").Value = 2936.1428
$ws.Range("L3").Value = 2936.1428
$ws.Range("N3").Value = -3164.1428

$ws = $wb.Worksheets.Item("BSM")  # row 22
$ws.Range("H22").Value = 298
$ws.Range("I22").Value = 298
$ws.Range("K22").Value = 298
$ws.Range("M22").Value = -125

$ws = $wb.Worksheets.Item("BSM")  # row 86
$ws.Range("H86").Value = 1420
$ws.Range("I86").Value = 833
$ws.Range("K86").Value = 833
$ws.Range("M86").Value = 290

$ws = $wb.Worksheets.Item("BSM")  # row 89
$ws.Range("H89").Value = 1420
$ws.Range("I89").Value = 833
$ws.Range("K89").Value = 4165
$ws.Range("M89").Value = 1451

$ws = $wb.Worksheets.Item("BSM")  # row 105
$ws.Range("H105").Value = 3006
$ws.Range("I105").Value = 2514.1
$ws.Range("K105").Value = 2514.1
$ws.Range("M105").Value = -767.0999999999999

$ws = $wb.Worksheets.Item("BSM")  # row 107
$ws.Range("H107").Value = 3989.5
$ws.Range("I107").Value = 3987.3076
$ws.Range("J107").Value = 3999
$ws.Range("K107").Value = 3987.3076
$ws.Range("L107").Value = 3999
$ws.Range("M107").Value = -2067.3076
$ws.Range("N107").Value = -7839

$ws = $wb.Worksheets.Item("CRP")  # row 21
$ws.Range("H21").Value = 4998.3335
$ws.Range("I21").Value = 4998.3335
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 4998.3335
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -4763.3335
$ws.Range("N21").ClearContents()

$ws = $wb.Worksheets.Item("CRP")  # row 31
$ws.Range("H31").Value = 3965.8909
$ws.Range("I31").Value = 1615.25
$ws.Range("K31").Value = 1615.25
$ws.Range("M31").Value = -1320.25

$ws = $wb.Worksheets.Item("CRP")  # row 34
$ws.Range("H34").Value = 3965.8909
$ws.Range("I34").Value = 1615.25
$ws.Range("K34").Value = 1615.25
$ws.Range("M34").Value = -1413.25

$ws = $wb.Worksheets.Item("CRP")  # row 94
$ws.Range("H94").Value = 1536.7646
$ws.Range("J94").Value = 1680.8572
$ws.Range("L94").Value = 1680.8572
$ws.Range("N94").Value = -2582.8572

$ws = $wb.Worksheets.Item("CUL")  # row 97
$ws.Range("H97").Value = 339.33334
$ws.Range("I97").Value = 466.5
$ws.Range("K97").Value = 1399.5
$ws.Range("M97").Value = -903.5

$ws = $wb.Worksheets.Item("GSM")  # row 46
$ws.Range("H46").Value = 20000
$ws.Range("J46").Value = 20000
$ws.Range("L46").Value = 20000
$ws.Range("N46").Value = -20312

$ws = $wb.Worksheets.Item("GSM")  # row 70
$ws.Range("H70").Value = 250332.83
$ws.Range("I70").Value = 299249.5
$ws.Range("J70").Value = 152499.5
$ws.Range("K70").Value = 299249.5
$ws.Range("L70").Value = 152499.5
$ws.Range("M70").Value = -298979.5
$ws.Range("N70").Value = -153039.5

$ws = $wb.Worksheets.Item("GSM")  # row 73
$ws.Range("H73").Value = 250332.83
$ws.Range("I73").Value = 299249.5
$ws.Range("J73").Value = 152499.5
$ws.Range("K73").Value = 299249.5
$ws.Range("L73").Value = 152499.5
$ws.Range("M73").Value = -298313.5
$ws.Range("N73").Value = -154371.5

$ws = $wb.Worksheets.Item("GSM")  # row 80
$ws.Range("H80").Value = 2484.3333
$ws.Range("I80").Value = 2225.9
$ws.Range("J80").Value = 3001.2
$ws.Range("K80").Value = 2225.9
$ws.Range("L80").Value = 3001.2
$ws.Range("M80").Value = -1227.9
$ws.Range("N80").Value = -4997.2

$ws = $wb.Worksheets.Item("GSM")  # row 83
$ws.Range("H83").Value = 2484.3333
$ws.Range("I83").Value = 2225.9
$ws.Range("J83").Value = 3001.2
$ws.Range("K83").Value = 11129.5
$ws.Range("L83").Value = 15006
$ws.Range("M83").Value = -6137.5
$ws.Range("N83").Value = -24990

$ws = $wb.Worksheets.Item("GSM")  # row 92
$ws.Range("H92").Value = 250251
$ws.Range("J92").Value = 250251
$ws.Range("L92").Value = 250251
$ws.Range("N92").Value = -253995

$ws = $wb.Worksheets.Item("GSM")  # row 102
$ws.Range("H102").Value = 1594.6666
$ws.Range("I102").Value = 1501.12
$ws.Range("K102").Value = 1501.12
$ws.Range("M102").Value = 120.8800000000001

$ws = $wb.Worksheets.Item("LTW")  # row 16
$ws.Range("H16").Value = 2373.2104
$ws.Range("I16").Value = 2240.647
$ws.Range("J16").Value = 3500
$ws.Range("K16").Value = 2240.647
$ws.Range("L16").Value = 3500
$ws.Range("M16").Value = -2070.647
$ws.Range("N16").Value = -3840

$ws = $wb.Worksheets.Item("LTW")  # row 68
$ws.Range("H68").Value = 1761.9286
$ws.Range("I68").Value = 1692.1818
$ws.Range("J68").Value = 2017.6666
$ws.Range("K68").Value = 1692.1818
$ws.Range("L68").Value = 2017.6666
$ws.Range("M68").Value = -943.1818000000001
$ws.Range("N68").Value = -3515.6666

$ws = $wb.Worksheets.Item("LTW")  # row 71
$ws.Range("H71").Value = 1761.9286
$ws.Range("I71").Value = 1692.1818
$ws.Range("J71").Value = 2017.6666
$ws.Range("K71").Value = 8460.909
$ws.Range("L71").Value = 10088.333
$ws.Range("M71").Value = -4716.909
$ws.Range("N71").Value = -17576.333

$ws = $wb.Worksheets.Item("LTW")  # row 93
$ws.Range("H93").Value = 1908.0588
$ws.Range("I93").Value = 1488.909
$ws.Range("K93").Value = 1488.909
$ws.Range("M93").Value = -240.9090000000001

$ws = $wb.Worksheets.Item("LTW")  # row 122
$ws.Range("H122").Value = 14229.412
$ws.Range("I122").Value = 13582.333
$ws.Range("J122").Value = 15782.4
$ws.Range("K122").Value = 40746.999
$ws.Range("L122").Value = 47347.2
$ws.Range("M122").Value = -38296.999
$ws.Range("N122").Value = -52247.2

$ws = $wb.Worksheets.Item("WVR")  # row 62
$ws.Range("H62").Value = 4445
$ws.Range("I62").Value = 3125
$ws.Range("J62").Value = 6425
$ws.Range("K62").Value = 3125
$ws.Range("L62").Value = 6425
$ws.Range("M62").Value = -2501
$ws.Range("N62").Value = -7673

$ws = $wb.Worksheets.Item("WVR")  # row 65
$ws.Range("H65").Value = 4445
$ws.Range("I65").Value = 3125
$ws.Range("J65").Value = 6425
$ws.Range("K65").Value = 15625
$ws.Range("L65").Value = 32125
$ws.Range("M65").Value = -12505
$ws.Range("N65").Value = -38365

$ws = $wb.Worksheets.Item("WVR")  # row 81
$ws.Range("H81").Value = 67102.81
$ws.Range("I81").Value = 203731.2
$ws.Range("J81").Value = 4999
$ws.Range("K81").Value = 407462.4
$ws.Range("L81").Value = 9998
$ws.Range("M81").Value = -406401.4
$ws.Range("N81").Value = -12120

$ws = $wb.Worksheets.Item("WVR")  # row 84
$ws.Range("H84").Value = 67102.81
$ws.Range("I84").Value = 203731.2
$ws.Range("J84").Value = 4999
$ws.Range("K84").Value = 2037312
$ws.Range("L84").Value = 49990
$ws.Range("M84").Value = -2032008
$ws.Range("N84").Value = -60598

$ws = $wb.Worksheets.Item("WVR")  # row 113
$ws.Range("H113").Value = 699.375
$ws.Range("I113").Value = 519.6
$ws.Range("K113").Value = 1558.8
$ws.Range("M113").Value = 611.1999999999998

